$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.63
$ws.Range("I2").Value = 2.8
$ws.Range("L2").Value = 3.5
$ws.Range("Z2").Value = 26
$ws.Range("AC2").Value = 8.5
$ws.Range("AG2").Value = 251
$ws.Range("AR2").Value = 81
$ws.Range("AU2").Value = 8
$ws.Range("AX2").Value = 15
$ws.Range("AY2").Value = 26
